$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 123; existing rows 123..177 shift down to 124..178.
$ws.Rows.Item(123).Insert()

# Populate the newly inserted row 123 with the new record's data.
$ws.Cells.Item(123, 1).Value = 10
$ws.Cells.Item(123, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(123, 3).Value = "La Araucanía"
$ws.Cells.Item(123, 4).Value = 45126
$ws.Cells.Item(123, 5).Value = 9
$ws.Cells.Item(123, 6).Value = "Fruta"
$ws.Cells.Item(123, 7).Value = 100107
$ws.Cells.Item(123, 8).Value = "Otros"
$ws.Cells.Item(123, 9).Value = 100107002
$ws.Cells.Item(123, 10).Value = "Chirimoya"
$ws.Cells.Item(123, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(123, 12).Value = "Primera"
$ws.Cells.Item(123, 13).Value = 80
$ws.Cells.Item(123, 14).Value = 5000
$ws.Cells.Item(123, 15).Value = 5000
$ws.Cells.Item(123, 16).Value = 5000
$ws.Cells.Item(123, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(123, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(123, 19).Value = 5000
$ws.Cells.Item(123, 20).Value = 1
